# Daily GitHub Actions refresh of the crypto symbol list (cryptos.xlsx / Sheet1).
# Column D ("Price") holds numeric-looking text (t="inlineStr" in the source file),
# so every Price write below is prefixed with a literal apostrophe to force
# Excel to keep it as text instead of silently coercing it to a Number cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price refresh for rows 2-9 (coins unchanged, only quotes moved) ---
$ws.Range("D2").Value  = "'244.93"     # BNB
$ws.Range("D4").Value  = "'5.194"      # HuobiToken
$ws.Range("D5").Value  = "'0.05749"    # Cronos
$ws.Range("D6").Value  = "'6.462"      # KuCoinToken
$ws.Range("D7").Value  = "'3.249"      # GateToken
$ws.Range("D8").Value  = "'0.8130"     # MXToken
$ws.Range("D9").Value  = "'0.8677"     # FTXToken

# --- Rows 10-18: a new coin ("One") was inserted at rank 9, pushing
#     WazirX..CoinExToken down by one row each; refresh every column ---
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01016"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1377"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.06923"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03208"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03013"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09319"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.809"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001543"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04719"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- Remaining scattered price refreshes (rows 19-50, coins unchanged) ---
$ws.Range("D19").Value = "'0.006230"
$ws.Range("D20").Value = "'0.001234"
$ws.Range("D21").Value = "'0.004093"
$ws.Range("D22").Value = "'0.00008691"
$ws.Range("D23").Value = "'3.576"
$ws.Range("D24").Value = "'2.153"
$ws.Range("D26").Value = "'0.1298"
$ws.Range("D27").Value = "'0.0002326"
$ws.Range("D40").Value = "'0.03709"
$ws.Range("D41").Value = "'0.006253"
$ws.Range("D42").Value = "'0.1049"
$ws.Range("D43").Value = "'0.002597"
$ws.Range("D44").Value = "'0.007088"
$ws.Range("D45").Value = "'0.00005261"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.4296"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").Value = "'0.002056"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.0001998"
